$wb = $excel.ActiveWorkbook

$commitSha = "623be9ca224ca149ccc629adf195eeced7f3526c"
$newGuid   = "7e9e826f-84fb-45dc-ae0f-755c458e070a"
$newFile   = "$newGuid.md"
$newHash   = "1481f6b6f15ce3071759f47acf51c21cab2b40b1"
$newUrl    = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/$commitSha/e2e/$newFile"

# ---------------------------------------------------------------
# Sheet "Overview": add row 3, mirroring row 2, for the new file
# ---------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$loOverview = $wsOverview.ListObjects.Item(1)

$wsOverview.Range("A2:G2").Copy($wsOverview.Range("A3:G3"))
$loOverview.Resize($wsOverview.Range("A1:G3"))

$wsOverview.Range("A3").Value = $newFile
$wsOverview.Range("B3").Value = "e2e\" + $newFile
$wsOverview.Range("C3").Value = ".md"
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-08-20 10:46:15"

$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), $newUrl, "", "", "e2e\" + $newFile)

# ---------------------------------------------------------------
# Sheet "zh-cn": add row 3, mirroring row 2, for the new file
# ---------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")
$loZh = $wsZh.ListObjects.Item(1)

$wsZh.Range("A2:P2").Copy($wsZh.Range("A3:P3"))
$loZh.Resize($wsZh.Range("A1:P3"))

$wsZh.Range("A3").Value = $newFile
$wsZh.Range("B3").Value = ".md"
$wsZh.Range("C3").Value = "Ready for handoff"
$wsZh.Range("D3").Value = "e2e"
$wsZh.Range("E3").Value = "ht"
$wsZh.Range("F3").Value = "'False"
$wsZh.Range("G3").Value = $newGuid + "." + $newHash + ".zh-cn.xlf"
$wsZh.Range("H3").Value = "2016-08-20 10:46:11"
$wsZh.Range("K3").Value = "0001-01-01 00:00:00"
$wsZh.Range("M3").Value = "'True"
$wsZh.Range("O3").Value = "'False"

$wsZh.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("K3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsZh.Hyperlinks.Add($wsZh.Range("A3"), $newUrl, "", "", $newFile)

# ---------------------------------------------------------------
# Sheet "de-de": add row 3, mirroring row 2, for the new file
# ---------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")
$loDe = $wsDe.ListObjects.Item(1)

$wsDe.Range("A2:P2").Copy($wsDe.Range("A3:P3"))
$loDe.Resize($wsDe.Range("A1:P3"))

$wsDe.Range("A3").Value = $newFile
$wsDe.Range("B3").Value = ".md"
$wsDe.Range("C3").Value = "Ready for handoff"
$wsDe.Range("D3").Value = "e2e"
$wsDe.Range("E3").Value = "ht"
$wsDe.Range("F3").Value = "'False"
$wsDe.Range("G3").Value = $newGuid + "." + $newHash + ".de-de.xlf"
$wsDe.Range("H3").Value = "2016-08-20 10:46:15"
$wsDe.Range("K3").Value = "0001-01-01 00:00:00"
$wsDe.Range("M3").Value = "'True"
$wsDe.Range("O3").Value = "'False"

$wsDe.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("K3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsDe.Hyperlinks.Add($wsDe.Range("A3"), $newUrl, "", "", $newFile)
